# Update the "Översikt ARBOGA" sheet (Avverkningsanmälningar):
#  1. Bump the "Förändrad" (column C) date from serial 45186 (2023-09-17) to
#     serial 45188 (2023-09-19) for every existing data row (2..261).
#  2. Row 261 picks up an explicit row height (matches the rest of the sheet).
#  3. Four brand-new notification rows (262..265) are appended for Sveaskog.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column C ("Förändrad") now reads serial 45188 for every existing row ---
$ws.Range("C2:C261").Value = 45188

# --- 2. Row 261 gets the same explicit 15pt row height as its neighbours ---
$ws.Rows.Item(261).RowHeight = 15

# --- 3. Append the four new rows reported for ARBOGA / Sveaskog ---
$newRows = @(
    @{ Row = 262; A = "A 43802-2023"; B = 45187; C = 45188; G = 1;   HasHeight = $true },
    @{ Row = 263; A = "A 43806-2023"; B = 45187; C = 45188; G = 0.8; HasHeight = $true },
    @{ Row = 264; A = "A 43813-2023"; B = 45187; C = 45188; G = 0.5; HasHeight = $true },
    @{ Row = 265; A = "A 43812-2023"; B = 45187; C = 45188; G = 1;   HasHeight = $false }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A                 # A: Beteckning
    $ws.Cells.Item($row, 2).Value = $r.B                 # B: Datum
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 3).Value = $r.C                 # C: Förändrad
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 4).Value = "VÄSTMANLANDS LÄN"    # D: Län
    $ws.Cells.Item($row, 5).Value = "ARBOGA"              # E: Kommun
    $ws.Cells.Item($row, 6).Value = "Sveaskog"            # F: Markägare
    $ws.Cells.Item($row, 7).Value = $r.G                  # G: Area (ha)
    $ws.Cells.Item($row, 8).Value = 0                     # H: Fridlysta
    $ws.Cells.Item($row, 9).Value = 0                     # I: Signalarter
    $ws.Cells.Item($row, 10).Value = 0                    # J: NT
    $ws.Cells.Item($row, 11).Value = 0                    # K: VU
    $ws.Cells.Item($row, 12).Value = 0                    # L: EN
    $ws.Cells.Item($row, 13).Value = 0                    # M: CR
    $ws.Cells.Item($row, 14).Value = 0                    # N: RE
    $ws.Cells.Item($row, 15).Value = 0                    # O: Rödlistade
    $ws.Cells.Item($row, 16).Value = 0                    # P: Hotade
    $ws.Cells.Item($row, 17).Value = 0                    # Q: Alla arter
    $ws.Cells.Item($row, 18).WrapText = $true             # R: Artnamn (blank, wrapped)

    if ($r.HasHeight) {
        $ws.Rows.Item($row).RowHeight = 15
    }
}
